# Revert "stationarity calculations, fk and radon plots"
# Restores the NormNMFRatio column to its position right after LogICARatio
# (shifting SubgraphStat/Emergence back one slot) and rewrites the affected
# LogICARatio / NormNMFRatio / SubgraphStat / Emergence data columns (X:AA)
# for every data row, plus a handful of recalculated AlgebraicConnect (J)
# and StatRat (V40) floating point values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: reflect the shared-string reordering ---
# Before: Y1=SubgraphStat, Z1=Emergence, AA1=NormNMFRatio
# After:  Y1=NormNMFRatio, Z1=SubgraphStat, AA1=Emergence (X1=LogICARatio unchanged)
$ws.Range("Y1").Value = "NormNMFRatio"
$ws.Range("Z1").Value = "SubgraphStat"
$ws.Range("AA1").Value = "Emergence"

# --- Data rows 3-45: recomputed stationarity / NMF / subgraph values ---
# Row 3
$ws.Range("J3").Value = 0.999999999999998
$ws.Range("X3").Value = 1.04887005367672
$ws.Range("Y3").Value = 1.47806084097062
$ws.Range("Z3").Value = 0.0130686392313725
$ws.Range("AA3").Value = 0.497251009589878
# Row 4
$ws.Range("X4").Value = 1.34588594833288
$ws.Range("Y4").Value = 0.89568857487168
$ws.Range("Z4").Value = 0.0157442604741007
$ws.Range("AA4").Value = -0.0179800374934191
# Row 5
$ws.Range("X5").Value = -0.853645217176353
$ws.Range("Y5").Value = 0.487905759876115
$ws.Range("Z5").Value = 0.00878987146533087
$ws.Range("AA5").Value = -0.144898012301384
# Row 6
$ws.Range("J6").Value = 1
$ws.Range("X6").Value = 0.614796557862919
$ws.Range("Y6").Value = 1.31832895468212
$ws.Range("Z6").Value = 0.0124307554377872
$ws.Range("AA6").Value = 0.263660888901169
# Row 7
$ws.Range("X7").Value = 0.270097535571253
$ws.Range("Y7").Value = 0.75853602126271
$ws.Range("Z7").Value = 0.0124307554377872
$ws.Range("AA7").Value = -0.206038554769815
# Row 8
$ws.Range("X8").Value = 1.22434314905709
$ws.Range("Y8").Value = 2.01279597217926
$ws.Range("Z8").Value = 0.0070118917343666
$ws.Range("AA8").Value = 0.186687152903642
# Row 9
$ws.Range("X9").Value = 1.78523177094415
$ws.Range("Y9").Value = 0.908884467951915
$ws.Range("Z9").Value = 0.0145159949751711
$ws.Range("AA9").Value = 0.47136943766907
# Row 10
$ws.Range("J10").Value = 0.238465591590029
$ws.Range("X10").Value = -1.00281087588735
$ws.Range("Y10").Value = 1.72134344587126
$ws.Range("Z10").Value = 0.0131379693689775
$ws.Range("AA10").Value = -1
# Row 11
$ws.Range("X11").Value = 0.766294317803198
$ws.Range("Y11").Value = 0.817432300754284
$ws.Range("Z11").Value = 0.0150361523392413
$ws.Range("AA11").Value = 1
# Row 12
$ws.Range("X12").Value = 0.576757883719376
$ws.Range("Y12").Value = 1.22070400001624
$ws.Range("Z12").Value = 0.0140866989190509
$ws.Range("AA12").Value = 0.207610035607276
# Row 13
$ws.Range("X13").Value = -1.18834813707523
$ws.Range("Y13").Value = 1.26845788294828
$ws.Range("Z13").Value = 0.0154534567266276
$ws.Range("AA13").Value = -0.504563175973437
# Row 14
$ws.Range("X14").Value = -0.306078378431768
$ws.Range("Y14").Value = 0.603924264107491
$ws.Range("Z14").Value = 0.0160828854606541
$ws.Range("AA14").Value = 0.508278477296733
# Row 15
$ws.Range("X15").Value = 1.38816303780057
$ws.Range("Y15").Value = 3.47512635591994
$ws.Range("Z15").Value = 0.0113929755478287
$ws.Range("AA15").Value = 0.0202201338050722
# Row 16
$ws.Range("X16").Value = -0.979063228217563
$ws.Range("Y16").Value = 0.648501521447253
$ws.Range("Z16").Value = 0.0145604426441847
$ws.Range("AA16").Value = -0.71971545076833
# Row 17
$ws.Range("X17").Value = -0.285897814894531
$ws.Range("Y17").Value = 0.986712034676278
$ws.Range("Z17").Value = 0.014825334134991
$ws.Range("AA17").Value = 0.655269943305426
# Row 18
$ws.Range("J18").Value = 0.0646489545249885
$ws.Range("X18").Value = -1.08268963872178
$ws.Range("Y18").Value = 1.15017241711094
$ws.Range("Z18").Value = 0.0153540765528727
$ws.Range("AA18").Value = -0.248977703341267
# Row 19
$ws.Range("X19").Value = 2.41117600821318
$ws.Range("Y19").Value = 1.06896439663969
$ws.Range("Z19").Value = 0.0158764726629214
$ws.Range("AA19").Value = 0.524839242191961
# Row 20
$ws.Range("X20").Value = 1.40976557261538
$ws.Range("Y20").Value = 0.907380390841924
$ws.Range("Z20").Value = 0.0141937749294579
$ws.Range("AA20").Value = -0.367113848097586
# Row 21
$ws.Range("X21").Value = 0.562431634653335
$ws.Range("Y21").Value = 0.993876622141327
$ws.Range("Z21").Value = 0.0145117123390032
$ws.Range("AA21").Value = 0.0495124643559382
# Row 22
$ws.Range("X22").Value = 0.439728819521156
$ws.Range("Y22").Value = 1.14969210902986
$ws.Range("Z22").Value = 0.014594047644061
$ws.Range("AA22").Value = -0.0658626360806904
# Row 23
$ws.Range("X23").Value = -1.41525850893781
$ws.Range("Y23").Value = 0.823466610119958
$ws.Range("Z23").Value = 0.014389861392556
$ws.Range("AA23").Value = -0.459535599239836
# Row 24
$ws.Range("X24").Value = -2.70338416780498
$ws.Range("Y24").Value = 6.36640412963768
$ws.Range("Z24").Value = 0.0161524545736765
$ws.Range("AA24").Value = 0.455185284677532
# Row 25
$ws.Range("X25").Value = 0.354501699724576
$ws.Range("Y25").Value = 0.158850872775147
$ws.Range("Z25").Value = 0.0153914881537371
$ws.Range("AA25").Value = -0.3333077950483
# Row 26
$ws.Range("X26").Value = -1.49678311568755
$ws.Range("Y26").Value = 2.31070909967728
$ws.Range("Z26").Value = 0.0157083458970645
$ws.Range("AA26").Value = 0.0221415726462084
# Row 27
$ws.Range("X27").Value = 1.38151543494553
$ws.Range("Y27").Value = 0.487947303386618
$ws.Range("Z27").Value = 0.0154573970921823
$ws.Range("AA27").Value = 0.273807624878081
# Row 28
$ws.Range("X28").Value = 0.154986937700053
$ws.Range("Y28").Value = 0.938555672609553
$ws.Range("Z28").Value = 0.0159509673976893
$ws.Range("AA28").Value = -0.329200170677304
# Row 29
$ws.Range("X29").Value = -0.934557070631527
$ws.Range("Y29").Value = 1.19335042512633
$ws.Range("Z29").Value = 0.0163705868910213
$ws.Range("AA29").Value = -0.264164328800264
# Row 30
$ws.Range("X30").Value = 0.139227600919719
$ws.Range("Y30").Value = 0.946094252235949
$ws.Range("Z30").Value = 0.0156018477241692
$ws.Range("AA30").Value = 0.362055877466574
# Row 31
$ws.Range("X31").Value = -0.828633894425627
$ws.Range("Y31").Value = 0.904090206526913
$ws.Range("Z31").Value = 0.0155941129759143
$ws.Range("AA31").Value = -0.067430633608431
# Row 32
$ws.Range("X32").Value = -2.5875680790163
$ws.Range("Y32").Value = 1.05616235065536
$ws.Range("Z32").Value = 0.0155906223665728
$ws.Range("AA32").Value = -0.590480977984785
# Row 33
$ws.Range("X33").Value = 0.913687368325867
$ws.Range("Y33").Value = 0.989825199424454
$ws.Range("Z33").Value = 0.0164239032884816
$ws.Range("AA33").Value = 0.763038460876312
# Row 34
$ws.Range("X34").Value = 1.09501346185881
$ws.Range("Y34").Value = 0.963061159874462
$ws.Range("Z34").Value = 0.0144777445925457
$ws.Range("AA34").Value = -0.0483531435235984
# Row 35
$ws.Range("X35").Value = 0.121835458367477
$ws.Range("Y35").Value = 0.968744140452917
$ws.Range("Z35").Value = 0.0139078463835676
$ws.Range("AA35").Value = -0.423828515486821
# Row 36
$ws.Range("J36").Value = 0.0618976942766494
$ws.Range("X36").Value = -1.15232183184174
$ws.Range("Y36").Value = 1.10770841969219
$ws.Range("Z36").Value = 0.0159379240949826
$ws.Range("AA36").Value = 0.281001650138757
# Row 37
$ws.Range("X37").Value = -0.801510330820907
$ws.Range("Y37").Value = 0.847050069997274
$ws.Range("Z37").Value = 0.0145325579534586
$ws.Range("AA37").Value = -0.557419284141504
# Row 38
$ws.Range("X38").Value = 0.566837583827508
$ws.Range("Y38").Value = 1.0703991782136
$ws.Range("Z38").Value = 0.0165550813212645
$ws.Range("AA38").Value = 0.455274433787522
# Row 39
$ws.Range("J39").Value = 0.0768814325141799
$ws.Range("X39").Value = 0.48768383003101
$ws.Range("Y39").Value = 1.06547540132603
$ws.Range("Z39").Value = 0.0174159070567434
$ws.Range("AA39").Value = 0.175574785184368
# Row 40
$ws.Range("V40").Value = 0.909879266755755
$ws.Range("X40").Value = -0.587404322954513
$ws.Range("Y40").Value = 0.939497303142611
$ws.Range("Z40").Value = 0.0146197522283737
$ws.Range("AA40").Value = 0.172025904277018
# Row 41
$ws.Range("J41").Value = 0.048675410971366
$ws.Range("X41").Value = -1.69455934132459
$ws.Range("Y41").Value = 1.02138039698852
$ws.Range("Z41").Value = 0.0157797226669743
$ws.Range("AA41").Value = -0.354922469674246
# Row 42
$ws.Range("X42").Value = -0.480404067142626
$ws.Range("Y42").Value = 0.67358089270154
$ws.Range("Z42").Value = 0.00944419553063511
$ws.Range("AA42").Value = 0.2746462865345
# Row 43
$ws.Range("X43").Value = -2.20627013906616
$ws.Range("Y43").Value = 0.316045913814248
$ws.Range("Z43").Value = 0.00642756742316938
$ws.Range("AA43").Value = -0.760420169006862
# Row 44
$ws.Range("J44").Value = 0.18121326453262
$ws.Range("X44").Value = 1.25358684555304
$ws.Range("Y44").Value = 1.88897245302525
$ws.Range("Z44").Value = 0.00861579736784824
$ws.Range("AA44").Value = 0.836037129462003
# Row 45
$ws.Range("X45").Value = -0.854212213859845
$ws.Range("Y45").Value = 0.692824659123163
$ws.Range("Z45").Value = 0.0133162403328807
$ws.Range("AA45").Value = -0.187871788619365
